$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2024.10.23 ruoqiang threshold update
# warning_thre (column C) value is folded into abnormal_thre (column B)
# for rows 2-3, and the now-redundant warning_thre cells are cleared.
$ws.Range("B2").Value = 75
$ws.Range("C2").Value = ""

$ws.Range("B3").Value = 75
$ws.Range("C3").Value = ""

$ws.Range("K11").Select()
